$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "LDLR"
$ws.Range("B14").Value = $true
$ws.Range("C14").Value = "often have pathogenic CNVs - complete genomic targeting"

$ws.Range("A1").Select() | Out-Null
